$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 141, pushing the existing data (old rows
# 141-229) down to rows 143-231.
$ws.Rows("141:142").Insert()

# Copy formatting (date style) for column D from the row that used to be 141
# (now at 143) onto the two new rows, so the new date cells keep the same
# number format as the rest of the column.
$ws.Range("D143").Copy()
$ws.Range("D141:D142").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the two new data rows with the values from the commit.
$ws.Range("A141").Value = 9
$ws.Range("B141").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C141").Value = "Metropolitana"
$ws.Range("D141").Value = 44603
$ws.Range("E141").Value = 13
$ws.Range("F141").Value = 100112030
$ws.Range("G141").Value = "Poroto granado"
$ws.Range("H141").Value = "Sin especificar"
$ws.Range("I141").Value = "Primera"
$ws.Range("J141").Value = 61
$ws.Range("K141").Value = 23000
$ws.Range("L141").Value = 25000
$ws.Range("M141").Value = 23984
$ws.Range("N141").Value = "$/saco 25 kilos"
$ws.Range("O141").Value = "Región Metropolitana"
$ws.Range("P141").Value = 959
$ws.Range("Q141").Value = 25
$ws.Range("R141").Value = "Hortaliza"

$ws.Range("A142").Value = 9
$ws.Range("B142").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C142").Value = "Metropolitana"
$ws.Range("D142").Value = 44603
$ws.Range("E142").Value = 13
$ws.Range("F142").Value = 100112030
$ws.Range("G142").Value = "Poroto granado"
$ws.Range("H142").Value = "Sin especificar"
$ws.Range("I142").Value = "Primera"
$ws.Range("J142").Value = 79
$ws.Range("K142").Value = 22000
$ws.Range("L142").Value = 25000
$ws.Range("M142").Value = 23519
$ws.Range("N142").Value = "$/saco 25 kilos"
$ws.Range("O142").Value = "Región del Maule"
$ws.Range("P142").Value = 941
$ws.Range("Q142").Value = 25
$ws.Range("R142").Value = "Hortaliza"
